# Actualización 18 de Octubre
# Adds two new daily rows (2020-04-16 and 2020-04-17, serials 43937/43938)
# to both "Hoja1" (national totals by region) and "Hoja2" sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# New data rows: date serial, dia, then the 16 region columns (C..R) and total (S)
$hoja1Rows = @(
    @(43937, 45, 134, 73, 211, 13, 68, 345, 4915, 55, 223, 656, 578, 907, 154, 412, 7, 501, 9252),
    @(43938, 46, 142, 81, 226, 13, 68, 359, 5192, 55, 276, 667, 606, 944, 156, 416, 7, 522, 9730)
)

$hoja2Rows = @(
    @(43937, 45, 2, 0, 1, 0, 0, 4, 51, 0, 7, 11, 3, 22, 3, 6, 0, 6, 116),
    @(43938, 46, 2, 0, 1, 0, 0, 4, 58, 0, 8, 12, 3, 23, 3, 6, 0, 6, 126)
)

function Fill-Row($ws, $rowNum, $values) {
    $ws.Cells.Item($rowNum, 1).Value = $values[0]
    $ws.Cells.Item($rowNum, 1).NumberFormat = "DD/MM/YY"
    for ($col = 2; $col -le $values.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
    }
}

$r = 46
foreach ($row in $hoja1Rows) {
    Fill-Row $ws1 $r $row
    $r = $r + 1
}

$r = 46
foreach ($row in $hoja2Rows) {
    Fill-Row $ws2 $r $row
    $r = $r + 1
}

# Update selections / active sheet to match the post-edit view state.
# Hoja2 selection first, then Hoja1 last so Hoja1 stays the active tab.
$ws2.Activate()
$ws2.Range("A45").Select()

$ws1.Activate()
$ws1.Range("A51").Select()
